# Se actualizo la lista de paquetes
# Adds the "devtools" package to the inventory sheet (new cell A18),
# and moves the underline emphasis that used to sit on B17 ("Ckmeans.1d.dp")
# over to A17 ("reshape2"), matching the author's re-styling of row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 restyle -------------------------------------------------
# B17 currently holds the only cell using the "lone" style (s=3). Underlining
# it updates that style definition in place. Underlining A17 right after lets
# it pick up (dedupe onto) that same, now-underlined style. Finally dropping
# the underline back off B17 leaves it with the plain/default look, while
# A17 keeps the underline -- i.e. the underline moves from B17 to A17.
$ws.Range("B17").Font.Underline = $true
$ws.Range("A17").Font.Underline = $true
$ws.Range("B17").Font.Underline = $false

# --- New package row -------------------------------------------------
# "devtools" is a new package added to the list, placed in column A of the
# previously A-empty row 18 (which already had "car" in column B).
$ws.Range("A18").Value = "devtools"

# --- Selection ---------------------------------------------------------
# The author's last selected cell ends up being A17.
$ws.Range("A17").Select() | Out-Null
